# Swap the data (columns B:AC) between row 46 and row 47.
# Column A (the running id 44 / 45) stays in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng46 = $ws.Range("B46:AC46")
$rng47 = $ws.Range("B47:AC47")

$vals46 = $rng46.Value2
$vals47 = $rng47.Value2

$rng46.Value = $vals47
$rng47.Value = $vals46
